# Scheduled-runner refresh of market-price-derived columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Only columns H..N are touched; a handful of rows gain or lose their N column
# (LeveProfitHQ) cell entirely when the HQ price/profit becomes inapplicable/
# applicable, which is replicated below by assigning "" to clear a cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 837.5
$ws.Range("I4").Value = 616.5
$ws.Range("J4").Value = 1500.5
$ws.Range("K4").Value = 616.5
$ws.Range("L4").Value = 1500.5
$ws.Range("M4").Value = -502.5
$ws.Range("N4").Value = -1728.5

$ws.Range("H18").Value = 621.75
$ws.Range("I18").Value = 621.75
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 621.75
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -337.75
$ws.Range("N18").Value = ""

$ws.Range("H28").Value = 2680.375
$ws.Range("I28").Value = 2698.9312
$ws.Range("J28").Value = 2501
$ws.Range("K28").Value = 2698.9312
$ws.Range("L28").Value = 2501
$ws.Range("M28").Value = -2213.9312
$ws.Range("N28").Value = -3471

$ws.Range("H41").Value = 1086.9286
$ws.Range("I41").Value = 2262.8333
$ws.Range("J41").Value = 205
$ws.Range("K41").Value = 2262.8333
$ws.Range("L41").Value = 205
$ws.Range("M41").Value = -1822.8333
$ws.Range("N41").Value = -1085

$ws.Range("H58").Value = 762.5
$ws.Range("I58").Value = 625
$ws.Range("J58").Value = 1450
$ws.Range("K58").Value = 1875
$ws.Range("L58").Value = 4350
$ws.Range("M58").Value = -1725
$ws.Range("N58").Value = -4650

$ws.Range("H62").Value = 6374.8335
$ws.Range("I62").Value = 7833
$ws.Range("K62").Value = 7833
$ws.Range("M62").Value = -7209

$ws.Range("H65").Value = 6374.8335
$ws.Range("I65").Value = 7833
$ws.Range("K65").Value = 39165
$ws.Range("M65").Value = -36045

$ws.Range("H113").Value = 5949.0347
$ws.Range("J113").Value = 5442.778
$ws.Range("L113").Value = 5442.778
$ws.Range("N113").Value = -11950.778

$ws.Range("H125").Value = 4654.615
$ws.Range("I125").Value = 3775.111
$ws.Range("J125").Value = 6633.5
$ws.Range("K125").Value = 33975.999
$ws.Range("L125").Value = 59701.5
$ws.Range("M125").Value = -31515.999
$ws.Range("N125").Value = -64621.5

$ws.Range("H127").Value = 1549.9286
$ws.Range("J127").Value = 1774
$ws.Range("L127").Value = 5322
$ws.Range("N127").Value = -15242

$ws.Range("H129").Value = 3341.1667
$ws.Range("I129").Value = 2988
$ws.Range("J129").Value = 3458.889
$ws.Range("K129").Value = 8964
$ws.Range("L129").Value = 10376.667
$ws.Range("M129").Value = -3964
$ws.Range("N129").Value = -20376.667

$ws.Range("H138").Value = 3880.21
$ws.Range("I138").Value = 6127.75
$ws.Range("J138").Value = 3573.7273
$ws.Range("K138").Value = 18383.25
$ws.Range("L138").Value = 10721.1819
$ws.Range("M138").Value = -13243.25
$ws.Range("N138").Value = -21001.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 784.2917
$ws.Range("I2").Value = 741.9048
$ws.Range("J2").Value = 1081
$ws.Range("K2").Value = 741.9048
$ws.Range("L2").Value = 1081
$ws.Range("M2").Value = -628.9048
$ws.Range("N2").Value = -1307

$ws.Range("H32").Value = 10686.267
$ws.Range("I32").Value = 6913.9443
$ws.Range("J32").Value = 25775.555
$ws.Range("K32").Value = 6913.9443
$ws.Range("L32").Value = 25775.555
$ws.Range("M32").Value = -6626.9443
$ws.Range("N32").Value = -26349.555

$ws.Range("H45").Value = 2809.5833
$ws.Range("I45").Value = 1549.5
$ws.Range("J45").Value = 4069.6667
$ws.Range("K45").Value = 1549.5
$ws.Range("L45").Value = 4069.6667
$ws.Range("M45").Value = -1172.5
$ws.Range("N45").Value = -4823.6667

$ws.Range("H110").Value = 826.75
$ws.Range("I110").Value = 788.2727
$ws.Range("K110").Value = 788.2727
$ws.Range("M110").Value = 1256.7273

$ws.Range("H116").Value = 784.2917
$ws.Range("I116").Value = 741.9048
$ws.Range("J116").Value = 1081
$ws.Range("K116").Value = 741.9048
$ws.Range("L116").Value = 1081
$ws.Range("M116").Value = 1552.0952
$ws.Range("N116").Value = -5669

$ws.Range("H132").Value = 1815.5883
$ws.Range("J132").Value = 2766
$ws.Range("L132").Value = 8298
$ws.Range("N132").Value = -13358

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 784.2917
$ws.Range("I3").Value = 741.9048
$ws.Range("J3").Value = 1081
$ws.Range("K3").Value = 741.9048
$ws.Range("L3").Value = 1081
$ws.Range("M3").Value = -627.9048
$ws.Range("N3").Value = -1309

$ws.Range("H86").Value = 2219.1667
$ws.Range("I86").Value = 1949.8
$ws.Range("J86").Value = 2411.5715
$ws.Range("K86").Value = 1949.8
$ws.Range("L86").Value = 2411.5715
$ws.Range("M86").Value = -826.8
$ws.Range("N86").Value = -4657.5715

$ws.Range("H89").Value = 2219.1667
$ws.Range("I89").Value = 1949.8
$ws.Range("J89").Value = 2411.5715
$ws.Range("K89").Value = 9749
$ws.Range("L89").Value = 12057.8575
$ws.Range("M89").Value = -4133
$ws.Range("N89").Value = -23289.8575

$ws.Range("H134").Value = 6912.7734
$ws.Range("I134").Value = 6033.814
$ws.Range("K134").Value = 18101.442
$ws.Range("M134").Value = -15566.442

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 302.91177
$ws.Range("I7").Value = 318.70587
$ws.Range("J7").Value = 287.11765
$ws.Range("K7").Value = 318.70587
$ws.Range("L7").Value = 287.11765
$ws.Range("M7").Value = -205.70587
$ws.Range("N7").Value = -513.11765

$ws.Range("H22").Value = 935.4
$ws.Range("I22").Value = 641.4286
$ws.Range("K22").Value = 641.4286
$ws.Range("M22").Value = -291.4286

$ws.Range("H86").Value = 5446.9443
$ws.Range("I86").Value = 3373
$ws.Range("J86").Value = 6766.727
$ws.Range("K86").Value = 3373
$ws.Range("L86").Value = 6766.727
$ws.Range("M86").Value = -2250
$ws.Range("N86").Value = -9012.726999999999

$ws.Range("H89").Value = 5446.9443
$ws.Range("I89").Value = 3373
$ws.Range("J89").Value = 6766.727
$ws.Range("K89").Value = 16865
$ws.Range("L89").Value = 33833.635
$ws.Range("M89").Value = -11249
$ws.Range("N89").Value = -45065.635

$ws.Range("H107").Value = 310.5
$ws.Range("I107").Value = 310.5
$ws.Range("K107").Value = 310.5
$ws.Range("M107").Value = 1609.5

$ws.Range("H141").Value = 369041.8
$ws.Range("I141").Value = 65098.4
$ws.Range("J141").Value = 672985.2
$ws.Range("K141").Value = 65098.4
$ws.Range("L141").Value = 672985.2
$ws.Range("M141").Value = -59918.4
$ws.Range("N141").Value = -683345.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2067.5217
$ws.Range("J34").Value = 6916.3335
$ws.Range("L34").Value = 20749.0005
$ws.Range("N34").Value = -20917.0005

$ws.Range("H37").Value = 132083.03
$ws.Range("J37").Value = 132083.03
$ws.Range("L37").Value = 396249.09
$ws.Range("N37").Value = -396473.09

$ws.Range("H39").Value = 6067.4546
$ws.Range("I39").Value = 2595
$ws.Range("J39").Value = 6839.1113
$ws.Range("K39").Value = 7785
$ws.Range("L39").Value = 20517.3339
$ws.Range("M39").Value = -7491
$ws.Range("N39").Value = -21105.3339

$ws.Range("H107").Value = 981.2
$ws.Range("I107").Value = 393
$ws.Range("K107").Value = 1179
$ws.Range("M107").Value = 741

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 412.2
$ws.Range("I2").Value = 70.30768999999999
$ws.Range("J2").Value = 1047.1428
$ws.Range("K2").Value = 70.30768999999999
$ws.Range("L2").Value = 1047.1428
$ws.Range("M2").Value = 42.69231000000001
$ws.Range("N2").Value = -1273.1428

$ws.Range("H36").Value = 9643.75
$ws.Range("J36").Value = 10225.333
$ws.Range("L36").Value = 10225.333
$ws.Range("N36").Value = -11195.333

$ws.Range("H70").Value = 9292.294
$ws.Range("I70").Value = 12352.125
$ws.Range("K70").Value = 12352.125
$ws.Range("M70").Value = -12082.125

$ws.Range("H73").Value = 9292.294
$ws.Range("I73").Value = 12352.125
$ws.Range("K73").Value = 12352.125
$ws.Range("M73").Value = -11416.125

$ws.Range("H122").Value = 81762.69500000001
$ws.Range("I122").Value = 145742.42
$ws.Range("J122").Value = 7119.6665
$ws.Range("K122").Value = 437227.26
$ws.Range("L122").Value = 21358.9995
$ws.Range("M122").Value = -434777.26
$ws.Range("N122").Value = -26258.9995

$ws.Range("H132").Value = 3764.587
$ws.Range("I132").Value = 2449.9473
$ws.Range("K132").Value = 7349.841899999999
$ws.Range("M132").Value = -4819.841899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2251.2307
$ws.Range("I100").Value = 1143
$ws.Range("K100").Value = 1143
$ws.Range("M100").Value = -602

$ws.Range("H122").Value = 4455.737
$ws.Range("I122").Value = 4261.357
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 12784.071
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -10334.071
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 4095.0312
$ws.Range("I132").Value = 3846.75
$ws.Range("J132").Value = 4508.8335
$ws.Range("K132").Value = 11540.25
$ws.Range("L132").Value = 13526.5005
$ws.Range("M132").Value = -9010.25
$ws.Range("N132").Value = -18586.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -11108

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""

$ws.Range("H122").Value = 67516.95
$ws.Range("I122").Value = 74502.44500000001
$ws.Range("K122").Value = 223507.335
$ws.Range("M122").Value = -221057.335

$ws.Range("H132").Value = 3477.318
$ws.Range("I132").Value = 3433.3809
$ws.Range("K132").Value = 10300.1427
$ws.Range("M132").Value = -7770.1427

$ws.Range("H136").Value = 10903.448
$ws.Range("I136").Value = 12041.875
$ws.Range("J136").Value = 9502.308000000001
$ws.Range("K136").Value = 36125.625
$ws.Range("L136").Value = 28506.924
$ws.Range("M136").Value = -33575.625
$ws.Range("N136").Value = -33606.924
